$d = $word.ActiveDocument

# Locate the contact-info paragraph (email / linkedin / portfolio line) and
# insert a new, blank paragraph immediately after it - mirroring the author
# pressing Enter at the end of that line. The new paragraph inherits the
# same centered / zero-spacing formatting as its neighbours.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*gefenbar23@gmail.com*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $target.Range.InsertParagraphAfter()
}
